# Fix layout using Tailwindcss, commit v3
# Adds 6 new student rows (13-18) below the existing data (which ended at row 12),
# cloning the formatting of row 12 and filling in the per-row specific data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows (13-18) right after the last existing data row (12) ---
$ws.Range("A13:W18").Insert()

# --- Clone formatting (styles) from row 12 onto the new rows ---
$ws.Range("A12:W12").Copy()
$ws.Range("A13:W18").PasteSpecial(-4122)
$ws.Range("A13:A18").RowHeight = 33.6

# --- Columns that repeat the same value on every data row (copy straight from row 12) ---
$repeatCols = @(2,3,5,7,8,9,10,11,12,13,14,23)
for ($r = 13; $r -le 18; $r++) {
    foreach ($col in $repeatCols) {
        $ws.Cells.Item($r, $col).Value = $ws.Cells.Item(12, $col).Value2
    }
}

# --- Per-row specific data ---
# Student codes (column D) are filled in first for every row, then the names
# (column F), then the scores - this matches the order new entries were added
# to the shared-string table.
$ws.Cells.Item(13,4).Value = "01167`nANHU"
$ws.Cells.Item(14,4).Value = "01166`nANHU"
$ws.Cells.Item(15,4).Value = "01165`nANHU"
$ws.Cells.Item(16,4).Value = "01164`nANHU"
$ws.Cells.Item(17,4).Value = "01163`nANHU"
$ws.Cells.Item(18,4).Value = "01162`nANHU"

$ws.Cells.Item(14,6).Value = "Nguyễn Trần Tuấn"
$ws.Cells.Item(13,6).Value = "Phạm Thế Tuấn"
$ws.Cells.Item(16,6).Value = "Hào Thế Long"
$ws.Cells.Item(17,6).Value = "Lò Thế Khanh"
$ws.Cells.Item(18,6).Value = "Nguyễn Văn Tuấn"
$ws.Cells.Item(15,6).Value = "Nguyễn Thế Quân"

# Scores (O:S) and bonus (U)
$ws.Cells.Item(13,15).Value = 10
$ws.Cells.Item(13,16).Value = 20
$ws.Cells.Item(13,17).Value = 30
$ws.Cells.Item(13,18).Value = 30
$ws.Cells.Item(13,19).Value = 50
$ws.Cells.Item(13,21).Value = 10

$ws.Cells.Item(14,15).Value = 50
$ws.Cells.Item(14,16).Value = 50
$ws.Cells.Item(14,17).Value = 30
$ws.Cells.Item(14,18).Value = 30
$ws.Cells.Item(14,19).Value = 50
$ws.Cells.Item(14,21).Value = 10

$ws.Cells.Item(15,15).Value = 10
$ws.Cells.Item(15,16).Value = 20
$ws.Cells.Item(15,17).Value = 30
$ws.Cells.Item(15,18).Value = 30
$ws.Cells.Item(15,19).Value = 50

$ws.Cells.Item(16,15).Value = 10
$ws.Cells.Item(16,16).Value = 20
$ws.Cells.Item(16,17).Value = 40
$ws.Cells.Item(16,18).Value = 30
$ws.Cells.Item(16,19).Value = 50

$ws.Cells.Item(17,15).Value = 10
$ws.Cells.Item(17,16).Value = 20
$ws.Cells.Item(17,17).Value = 30
$ws.Cells.Item(17,18).Value = 50
$ws.Cells.Item(17,19).Value = 50
$ws.Cells.Item(17,21).Value = 10

$ws.Cells.Item(18,15).Value = 40
$ws.Cells.Item(18,16).Value = 20
$ws.Cells.Item(18,17).Value = 30
$ws.Cells.Item(18,18).Value = 30
$ws.Cells.Item(18,19).Value = 50
$ws.Cells.Item(18,21).Value = 10

# --- Formulas: STT (A), total score (T) and total+bonus (V) ---
for ($r = 13; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Formula = '=IF(B' + $r + '="","",SUBTOTAL(3,$B$6:$B' + $r + '))'
    $ws.Cells.Item($r, 20).Formula = "=SUM(O$r`:S$r)"
    $ws.Cells.Item($r, 22).Formula = "=T$r+U$r"
}

# --- Update the view: scroll so row 7 / column K is the top-left cell, select S15 ---
$ws.Application.Goto($ws.Range("K7"), $true)
$ws.Range("S15").Select()
